# The commit swaps the contents of ppt/theme/theme1.xml and
# ppt/theme/theme2.xml: the deck's live theme (the one actually wired to
# the slide master / presentation, i.e. theme2.xml) changes from the
# "Integral" / "Red Violet" color scheme to the plain "Office Theme" /
# "Office" color scheme (and vice-versa for the otherwise-unused
# notes-master theme part). Font scheme and format scheme are identical
# between the two theme parts already, so the only user-visible effect
# is the 12 theme colors used throughout every slide/layout/master.
#
# Reproduce that with the Design ColorScheme COM surface, which is the
# supported way to repaint a presentation's live theme colors.

$p = $ppt.ActivePresentation
$cs = $p.SlideMaster.ColorScheme

# ppColorSchemeIndex order: dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink.
# Target values are the "Office" scheme (was previously in theme1.xml).
$cs.Colors(1).RGB  = 0         # dk1      000000
$cs.Colors(2).RGB  = 16777215  # lt1      FFFFFF
$cs.Colors(3).RGB  = 6968388   # dk2      44546A
$cs.Colors(4).RGB  = 15132391  # lt2      E7E6E6
$cs.Colors(5).RGB  = 13998939  # accent1  5B9BD5
$cs.Colors(6).RGB  = 3243501   # accent2  ED7D31
$cs.Colors(7).RGB  = 10855845  # accent3  A5A5A5
$cs.Colors(8).RGB  = 49407     # accent4  FFC000
$cs.Colors(9).RGB  = 12874308  # accent5  4472C4
$cs.Colors(10).RGB = 4697456   # accent6  70AD47
$cs.Colors(11).RGB = 12673797  # hlink    0563C1
$cs.Colors(12).RGB = 7491477   # folHlink 954F72
